$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 509, shifting the existing rows 509:624 down to 510:625
$ws.Rows.Item(509).Insert()

# Populate the new row 509 with the new weekly record
$ws.Range("A509").Value = 5
$ws.Range("B509").Value = "Macroferia Regional de Talca"
$ws.Range("C509").Value = "Maule"
$ws.Range("D509").Value = 45173
$ws.Range("E509").Value = 7
$ws.Range("F509").Value = 100112023
$ws.Range("G509").Value = "Brócoli"
$ws.Range("H509").Value = "Sin especificar"
$ws.Range("I509").Value = "Primera"
$ws.Range("J509").Value = 3000
$ws.Range("K509").Value = 700
$ws.Range("L509").Value = 700
$ws.Range("M509").Value = 700
$ws.Range("N509").Value = "$/unidad"
$ws.Range("O509").Value = "Región del Maule"
$ws.Range("P509").Value = 700
$ws.Range("Q509").Value = 1
$ws.Range("R509").Value = "Hortaliza"
